$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计"
$ws2 = $wb.Worksheets.Item(2)   # currently "2022-Q3", will become "2022-Q4"

# ------------------------------------------------------------------
# 1) Insert the new worksheet that will hold the (unchanged) old
#    "2022-Q3" fund-holding data, placed right after $ws2 so it lands
#    at position 3 / rId3 while $ws2 keeps rId2.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)

# Copy the current header row & column-A formatting (style "1") from
# $ws2 over to $ws3 before we overwrite $ws2 with the new Q4 data.
$ws2.Range("B1:H1").Copy($ws3.Range("B1:H1"))
$ws2.Range("A2:A6").Copy($ws3.Range("A2:A6"))

# Fill in the rest of the (unchanged) old 2022-Q3 fund data on $ws3.
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = '590001'
$ws3.Range("C2").Value = '中邮核心优选混合'
$ws3.Range("D2").NumberFormat = "@"
$ws3.Range("D2").Value = '12.91'
$ws3.Range("E2").NumberFormat = "@"
$ws3.Range("E2").Value = '81.83'
$ws3.Range("F2").NumberFormat = "@"
$ws3.Range("F2").Value = '3.04'
$ws3.Range("G2").NumberFormat = "@"
$ws3.Range("G2").Value = '0.3925'
$ws3.Range("H2").Value = 10
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = '590005'
$ws3.Range("C3").Value = '中邮核心主题混合'
$ws3.Range("D3").NumberFormat = "@"
$ws3.Range("D3").Value = '5.69'
$ws3.Range("E3").NumberFormat = "@"
$ws3.Range("E3").Value = '82.29'
$ws3.Range("F3").NumberFormat = "@"
$ws3.Range("F3").Value = '3.88'
$ws3.Range("G3").NumberFormat = "@"
$ws3.Range("G3").Value = '0.2208'
$ws3.Range("H3").Value = 5
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("B4").Value = '003842'
$ws3.Range("C4").Value = '中邮景泰灵活配置混合A'
$ws3.Range("D4").NumberFormat = "@"
$ws3.Range("D4").Value = '2.09'
$ws3.Range("E4").NumberFormat = "@"
$ws3.Range("E4").Value = '35.39'
$ws3.Range("F4").NumberFormat = "@"
$ws3.Range("F4").Value = '1.60'
$ws3.Range("G4").NumberFormat = "@"
$ws3.Range("G4").Value = '0.0334'
$ws3.Range("H4").Value = 4
$ws3.Range("B5").NumberFormat = "@"
$ws3.Range("B5").Value = '001430'
$ws3.Range("C5").Value = '中邮乐享收益灵活配置混合'
$ws3.Range("D5").NumberFormat = "@"
$ws3.Range("D5").Value = '0.13'
$ws3.Range("E5").NumberFormat = "@"
$ws3.Range("E5").Value = '29.80'
$ws3.Range("F5").NumberFormat = "@"
$ws3.Range("F5").Value = '1.39'
$ws3.Range("G5").NumberFormat = "@"
$ws3.Range("G5").Value = '0.0018'
$ws3.Range("H5").Value = 2
$ws3.Range("B6").NumberFormat = "@"
$ws3.Range("B6").Value = '003843'
$ws3.Range("C6").Value = '中邮景泰灵活配置混合C'
$ws3.Range("D6").NumberFormat = "@"
$ws3.Range("D6").Value = '0.06'
$ws3.Range("E6").NumberFormat = "@"
$ws3.Range("E6").Value = '35.39'
$ws3.Range("F6").NumberFormat = "@"
$ws3.Range("F6").Value = '1.60'
$ws3.Range("G6").NumberFormat = "@"
$ws3.Range("G6").Value = '0.0010'
$ws3.Range("H6").Value = 4

# ------------------------------------------------------------------
# 2) Re-style $ws2's header row & column A to style "2" (matching the
#    "总计" sheet's header style) and overwrite with the new 2022-Q4 data.
# ------------------------------------------------------------------
$ws1.Range("B1").Copy($ws2.Range("B1:H1"))
$ws1.Range("A2").Copy($ws2.Range("A2:A15"))

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = '590002'
$ws2.Range("C2").Value = '中邮核心成长混合'
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = '35.53'
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = '73.54'
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = '2.63'
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = '0.9344'
$ws2.Range("H2").Value = 9
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = '590001'
$ws2.Range("C3").Value = '中邮核心优选混合'
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = '11.94'
$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = '83.93'
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = '4.75'
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Value = '0.5672'
$ws2.Range("H3").Value = 4
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = '590005'
$ws2.Range("C4").Value = '中邮核心主题混合'
$ws2.Range("D4").NumberFormat = "@"
$ws2.Range("D4").Value = '5.51'
$ws2.Range("E4").NumberFormat = "@"
$ws2.Range("E4").Value = '89.10'
$ws2.Range("F4").NumberFormat = "@"
$ws2.Range("F4").Value = '4.52'
$ws2.Range("G4").NumberFormat = "@"
$ws2.Range("G4").Value = '0.2491'
$ws2.Range("H4").Value = 5
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = '003887'
$ws2.Range("C5").Value = '汇安丰利灵活配置混合C'
$ws2.Range("D5").NumberFormat = "@"
$ws2.Range("D5").Value = '2.74'
$ws2.Range("E5").NumberFormat = "@"
$ws2.Range("E5").Value = '94.57'
$ws2.Range("F5").NumberFormat = "@"
$ws2.Range("F5").Value = '4.93'
$ws2.Range("G5").NumberFormat = "@"
$ws2.Range("G5").Value = '0.1351'
$ws2.Range("H5").Value = 4
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = '015092'
$ws2.Range("C6").Value = '汇安远见成长混合A'
$ws2.Range("D6").NumberFormat = "@"
$ws2.Range("D6").Value = '1.88'
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = '84.66'
$ws2.Range("F6").NumberFormat = "@"
$ws2.Range("F6").Value = '4.77'
$ws2.Range("G6").NumberFormat = "@"
$ws2.Range("G6").Value = '0.0897'
$ws2.Range("H6").Value = 4
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = '003886'
$ws2.Range("C7").Value = '汇安丰利灵活配置混合A'
$ws2.Range("D7").NumberFormat = "@"
$ws2.Range("D7").Value = '1.44'
$ws2.Range("E7").NumberFormat = "@"
$ws2.Range("E7").Value = '94.57'
$ws2.Range("F7").NumberFormat = "@"
$ws2.Range("F7").Value = '4.93'
$ws2.Range("G7").NumberFormat = "@"
$ws2.Range("G7").Value = '0.0710'
$ws2.Range("H7").Value = 4
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = '003842'
$ws2.Range("C8").Value = '中邮景泰灵活配置混合A'
$ws2.Range("D8").NumberFormat = "@"
$ws2.Range("D8").Value = '1.39'
$ws2.Range("E8").NumberFormat = "@"
$ws2.Range("E8").Value = '37.03'
$ws2.Range("F8").NumberFormat = "@"
$ws2.Range("F8").Value = '2.60'
$ws2.Range("G8").NumberFormat = "@"
$ws2.Range("G8").Value = '0.0361'
$ws2.Range("H8").Value = 1
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = '007775'
$ws2.Range("C9").Value = '汇安量化先锋混合A'
$ws2.Range("D9").NumberFormat = "@"
$ws2.Range("D9").Value = '0.22'
$ws2.Range("E9").NumberFormat = "@"
$ws2.Range("E9").Value = '86.70'
$ws2.Range("F9").NumberFormat = "@"
$ws2.Range("F9").Value = '3.83'
$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = '0.0084'
$ws2.Range("H9").Value = 6
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = '007776'
$ws2.Range("C10").Value = '汇安量化先锋混合C'
$ws2.Range("D10").NumberFormat = "@"
$ws2.Range("D10").Value = '0.12'
$ws2.Range("E10").NumberFormat = "@"
$ws2.Range("E10").Value = '86.70'
$ws2.Range("F10").NumberFormat = "@"
$ws2.Range("F10").Value = '3.83'
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = '0.0046'
$ws2.Range("H10").Value = 6
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = '003855'
$ws2.Range("C11").Value = '汇安丰华灵活配置混合C'
$ws2.Range("D11").NumberFormat = "@"
$ws2.Range("D11").Value = '0.18'
$ws2.Range("E11").NumberFormat = "@"
$ws2.Range("E11").Value = '45.67'
$ws2.Range("F11").NumberFormat = "@"
$ws2.Range("F11").Value = '2.35'
$ws2.Range("G11").NumberFormat = "@"
$ws2.Range("G11").Value = '0.0042'
$ws2.Range("H11").Value = 8
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = '001430'
$ws2.Range("C12").Value = '中邮乐享收益灵活配置混合'
$ws2.Range("D12").NumberFormat = "@"
$ws2.Range("D12").Value = '0.14'
$ws2.Range("E12").NumberFormat = "@"
$ws2.Range("E12").Value = '27.17'
$ws2.Range("F12").NumberFormat = "@"
$ws2.Range("F12").Value = '2.43'
$ws2.Range("G12").NumberFormat = "@"
$ws2.Range("G12").Value = '0.0034'
$ws2.Range("H12").Value = 1
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = '015093'
$ws2.Range("C13").Value = '汇安远见成长混合C'
$ws2.Range("D13").NumberFormat = "@"
$ws2.Range("D13").Value = '0.05'
$ws2.Range("E13").NumberFormat = "@"
$ws2.Range("E13").Value = '84.66'
$ws2.Range("F13").NumberFormat = "@"
$ws2.Range("F13").Value = '4.77'
$ws2.Range("G13").NumberFormat = "@"
$ws2.Range("G13").Value = '0.0024'
$ws2.Range("H13").Value = 4
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = '003843'
$ws2.Range("C14").Value = '中邮景泰灵活配置混合C'
$ws2.Range("D14").NumberFormat = "@"
$ws2.Range("D14").Value = '0.05'
$ws2.Range("E14").NumberFormat = "@"
$ws2.Range("E14").Value = '37.03'
$ws2.Range("F14").NumberFormat = "@"
$ws2.Range("F14").Value = '2.60'
$ws2.Range("G14").NumberFormat = "@"
$ws2.Range("G14").Value = '0.0013'
$ws2.Range("H14").Value = 1
$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = '003854'
$ws2.Range("C15").Value = '汇安丰华灵活配置混合A'
$ws2.Range("D15").NumberFormat = "@"
$ws2.Range("D15").Value = '0.00'
$ws2.Range("E15").NumberFormat = "@"
$ws2.Range("E15").Value = '45.67'
$ws2.Range("F15").NumberFormat = "@"
$ws2.Range("F15").Value = '2.35'
$ws2.Range("G15").Value = 0
$ws2.Range("H15").Value = 8

# Column-A row indices (0-based) for the new Q4 data
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("A5").Value = 3
$ws2.Range("A6").Value = 4
$ws2.Range("A7").Value = 5
$ws2.Range("A8").Value = 6
$ws2.Range("A9").Value = 7
$ws2.Range("A10").Value = 8
$ws2.Range("A11").Value = 9
$ws2.Range("A12").Value = 10
$ws2.Range("A13").Value = 11
$ws2.Range("A14").Value = 12
$ws2.Range("A15").Value = 13

# ------------------------------------------------------------------
# 3) Rename the sheets to their final names.
# ------------------------------------------------------------------
$ws2.Name = "2022-Q4"
$ws3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 4) Update the "总计" (summary) sheet: row 2 becomes the new 2022-Q4
#    totals, and a new row 3 is added with the old 2022-Q3 totals
#    (previously stored in row 2).
# ------------------------------------------------------------------
$ws1.Range("A2").Copy($ws1.Range("A3"))
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q3"
$ws1.Range("C3").Value = 5
$ws1.Range("D3").Value = 0.65

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 14
$ws1.Range("D2").Value = 2.11
